# Swap batting/bowling scorecards: England innings moved from columns A-F
# to J-O and Sri Lanka innings moved from J-O to A-F (the two team blocks
# were recorded under the wrong team after the tournament closure fix),
# together with the corrected runs/balls/dismissal figures for each player.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: opening batters
$ws.Range("A2").Value = 'Dinesh Chandimal'
$ws.Range("B2").Value = 19
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 'LBW'
$ws.Range("E2").Value = ' Chris Jordan'
$ws.Range("J2").Value = 'Jason Roy'
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 'Bowled'
$ws.Range("N2").Value = ' Chamika Karunarathne'

# Row 3
$ws.Range("A3").Value = 'Pathum Nissanka'
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("E3").Value = ' Mark Wood'
$ws.Range("J3").Value = 'Jos Buttler'
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 'Caught'
$ws.Range("N3").Value = ' Maheesh Theekshana'

# Row 4
$ws.Range("A4").Value = 'Charith Asalanka'
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 'Caught'
$ws.Range("E4").Value = ' Mark Wood'
$ws.Range("J4").Value = 'Dawid Malan'
$ws.Range("K4").Value = 38
$ws.Range("L4").Value = 11
$ws.Range("N4").Value = ' Chamika Karunarathne'

# Row 5
$ws.Range("A5").Value = 'Dhananjaya de Silva'
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 'Bowled'
$ws.Range("E5").Value = ' Adil Rashid'
$ws.Range("J5").Value = 'Jonny Bairstow'
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 'LBW'
$ws.Range("N5").Value = ' Chamika Karunarathne'

# Row 6
$ws.Range("A6").Value = 'Bhanuka Rajapakse'
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 'Caught'
$ws.Range("E6").Value = ' Chris Woakes'
$ws.Range("J6").Value = 'Eoin Morgan(C)'
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 'Caught'
$ws.Range("N6").Value = ' Wanindu Hasaranga'

# Row 7
$ws.Range("A7").Value = 'Dasun Shanka(C)'
$ws.Range("B7").Value = 16
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 'Caught'
$ws.Range("E7").Value = ' Adil Rashid'
$ws.Range("J7").Value = 'Moeen Ali'
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1
$ws.Range("N7").Value = ' Wanindu Hasaranga'

# Row 8
$ws.Range("A8").Value = 'Wanindu Hasaranga'
$ws.Range("B8").Value = 20
$ws.Range("C8").Value = 7
$ws.Range("E8").Value = ' Liam Livingstone'
$ws.Range("J8").Value = 'Liam Livingstone'
$ws.Range("K8").Value = 16
$ws.Range("L8").Value = 6
$ws.Range("M8").Value = 'NOT OUT'
$ws.Range("N8").Value = ' '

# Row 9
$ws.Range("A9").Value = 'Chamika Karunarathne'
$ws.Range("B9").Value = 16
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 'Bowled'
$ws.Range("E9").Value = ' Chris Jordan'
$ws.Range("J9").Value = 'Chris Woakes'
$ws.Range("K9").Value = 9
$ws.Range("M9").Value = '* NOT OUT'
$ws.Range("N9").Value = ' '

# Row 10
$ws.Range("A10").Value = 'Dushmantha Chameera'
$ws.Range("B10").Value = 2
$ws.Range("D10").Value = 'Caught'
$ws.Range("E10").Value = ' Adil Rashid'
$ws.Range("J10").Value = 'Chris Jordan'
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = ' '

# Row 11
$ws.Range("A11").Value = 'Maheesh Theekshana'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 'Bowled'
$ws.Range("E11").Value = ' Chris Jordan'
$ws.Range("J11").Value = 'Adil Rashid'

# Row 12
$ws.Range("A12").Value = 'Nuwan Pradeep'
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 2
$ws.Range("J12").Value = 'Mark Wood'

# Innings totals
$ws.Range("A16").Value = 86
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '8.0'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 48
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = 6
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = '5.4'
$ws.Range("L16").Style = "Normal"
$ws.Range("M16").Value = 34

# Bowling row 21
$ws.Range("A21").Value = 'Chris Woakes'
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 6
$ws.Range("J21").Value = 'Maheesh Theekshana'
$ws.Range("L21").Value = 16
$ws.Range("N21").Value = 16

# Bowling row 22
$ws.Range("A22").Value = 'Liam Livingstone'
$ws.Range("C22").Value = 14
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 14
$ws.Range("J22").Value = 'Dushmantha Chameera'
$ws.Range("L22").Value = 24
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 24

# Bowling row 23
$ws.Range("A23").Value = 'Mark Wood'
$ws.Range("C23").Value = 24
$ws.Range("E23").Value = 12
$ws.Range("J23").Value = 'Chamika Karunarathne'
$ws.Range("L23").Value = 6
$ws.Range("M23").Value = 3
$ws.Range("N23").Value = 6

# Bowling row 24
$ws.Range("A24").Value = 'Adil Rashid'
$ws.Range("J24").Value = 'Wanindu Hasaranga'
$ws.Range("L24").Value = 14
$ws.Range("M24").Value = 2
$ws.Range("N24").Value = 14

# Bowling row 25
$ws.Range("A25").Value = 'Chris Jordan'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '2.0'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 9.5
$ws.Range("J25").Value = 'Nuwan Pradeep'
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = '1.4'
$ws.Range("K25").Style = "Normal"
$ws.Range("L25").Value = 29
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 20.71
